# Regenerate merged AHB files: rename the "_old"/"_new" header suffixes to
# the concrete format versions being compared ("_FV2404"/"_FV2410"), freeze
# the header row, and turn the data range into a proper Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# Column base-names, in the order they appear for both the "_old"/"_FV2404"
# block (A1:J1) and the "_new"/"_FV2410" block (L1:U1). Column K1 is "diff"
# and is left untouched.
$suffixBases = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $suffixBases.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $suffixBases[$i] + "_FV2404"   # A1:J1
}
for ($i = 0; $i -lt $suffixBases.Count; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $suffixBases[$i] + "_FV2410"  # L1:U1
}

# Freeze panes so the header row (row 1) stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Convert the full used range into an Excel Table ("Table1") with an
# auto filter on the header row, matching the regenerated workbook layout.
$dataRange = $ws.Range("A1:U56")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

Write-Host "Edit complete"
